$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the simple appointment IDs (A00x -> AP0x) for the first three
# data rows, which keep their row position. ---
$ws.Range("A2").Value = "AP01"
$ws.Range("A3").Value = "AP02"
$ws.Range("A4").Value = "AP03"

# --- Row 5 and row 6 swap their data: the appointment that used to sit in
# row 6 (AP05 / "COMPLETED") now shows in row 5, and the appointment that
# used to sit in row 5 (AP04 / "SCHEDULED") now shows in row 6 with an
# updated date/time. Clear both rows first so no stale number formatting
# is left behind on cells that no longer need it. ---
$ws.Range("A5:G6").Clear()

# New row 5 = old row 6 data, renamed A005 -> AP05
$ws.Range("A5").Value = "AP05"
$ws.Range("B5").Value = "P1001"
$ws.Range("C5").Value = "D002"
$ws.Range("D5").Value = "COMPLETED"
$ws.Range("E5").Value = [DateTime]"2024-10-30"
$ws.Range("E5").NumberFormat = "mm-dd-yy"
$ws.Range("F5").Value = 0.41666666666666669
$ws.Range("F5").NumberFormat = "h:mm:ss AM/PM"
$ws.Range("G5").Value = "Health check up, all normal"

# New row 6 = old row 5 data, renamed A004 -> AP04, appointment moved to a
# new date/time (2024-10-31, 2:00 pm)
$ws.Range("A6").Value = "AP04"
$ws.Range("B6").Value = "P1001"
$ws.Range("C6").Value = "D001"
$ws.Range("D6").Value = "SCHEDULED"
$ws.Range("E6").Value = 45596
$ws.Range("F6").Value = "2:00 pm"
$ws.Range("G6").Value = ""

# Touch column H (no visible value) so the sheet's used range grows to
# include it, matching the workbook's recorded dimension.
$ws.Range("H6").NumberFormat = "General"

# Restore the cursor/selection position left in the file.
$ws.Range("C9").Select()
